$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 141.33333
$ws.Cells.Item(12, 9).Value = 149.6
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 149.6
$ws.Cells.Item(12, 12).Value = 100
$ws.Cells.Item(12, 13).Value = 20.40000000000001
$ws.Cells.Item(12, 14).Value = -440

$ws.Cells.Item(33, 8).Value = 248.81818
$ws.Cells.Item(33, 9).Value = 136.125
$ws.Cells.Item(33, 10).Value = 549.3333
$ws.Cells.Item(33, 11).Value = 136.125
$ws.Cells.Item(33, 12).Value = 549.3333
$ws.Cells.Item(33, 13).Value = 92.875
$ws.Cells.Item(33, 14).Value = -1007.3333

$ws.Cells.Item(64, 8).Value = 4590
$ws.Cells.Item(64, 9).Value = 5495
$ws.Cells.Item(64, 10).Value = 3986.6667
$ws.Cells.Item(64, 11).Value = 5495
$ws.Cells.Item(64, 12).Value = 3986.6667
$ws.Cells.Item(64, 13).Value = -5247
$ws.Cells.Item(64, 14).Value = -4482.6667

$ws.Cells.Item(67, 8).Value = 4590
$ws.Cells.Item(67, 9).Value = 5495
$ws.Cells.Item(67, 10).Value = 3986.6667
$ws.Cells.Item(67, 11).Value = 5495
$ws.Cells.Item(67, 12).Value = 3986.6667
$ws.Cells.Item(67, 13).Value = -4637
$ws.Cells.Item(67, 14).Value = -5702.6667

$ws.Cells.Item(97, 8).Value = 2273.6
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 2273.6
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 6820.799999999999
$ws.Cells.Item(97, 14).Value = -7812.799999999999

$ws.Cells.Item(112, 8).Value = 2458.862
$ws.Cells.Item(112, 9).Value = 749.8333
$ws.Cells.Item(112, 10).Value = 2904.6956
$ws.Cells.Item(112, 11).Value = 2249.4999
$ws.Cells.Item(112, 12).Value = 8714.086800000001
$ws.Cells.Item(112, 13).Value = -1141.4999
$ws.Cells.Item(112, 14).Value = -10930.0868

$ws.Cells.Item(118, 8).Value = 1069.7778
$ws.Cells.Item(118, 9).Value = 519.7143
$ws.Cells.Item(118, 10).Value = 2995
$ws.Cells.Item(118, 11).Value = 1559.1429
$ws.Cells.Item(118, 12).Value = 8985
$ws.Cells.Item(118, 13).Value = 97.85710000000017
$ws.Cells.Item(118, 14).Value = -12299

$ws.Cells.Item(132, 8).Value = 7097797.5
$ws.Cells.Item(132, 9).Value = 9808872
$ws.Cells.Item(132, 10).Value = 7294.077
$ws.Cells.Item(132, 11).Value = 29426616
$ws.Cells.Item(132, 12).Value = 21882.231
$ws.Cells.Item(132, 13).Value = -29424086
$ws.Cells.Item(132, 14).Value = -26942.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3331.7463
$ws.Cells.Item(32, 9).Value = 2980.1094
$ws.Cells.Item(32, 10).Value = 10833.333
$ws.Cells.Item(32, 11).Value = 2980.1094
$ws.Cells.Item(32, 12).Value = 10833.333
$ws.Cells.Item(32, 13).Value = -2693.1094
$ws.Cells.Item(32, 14).Value = -11407.333

$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 14).ClearContents()

$ws.Cells.Item(61, 8).Value = 1250.579
$ws.Cells.Item(61, 9).Value = 1157.3715
$ws.Cells.Item(61, 10).Value = 2338
$ws.Cells.Item(61, 11).Value = 1157.3715
$ws.Cells.Item(61, 12).Value = 2338
$ws.Cells.Item(61, 13).Value = -945.3715
$ws.Cells.Item(61, 14).Value = -2762

$ws.Cells.Item(132, 8).Value = 1316.2407
$ws.Cells.Item(132, 9).Value = 1045.1163
$ws.Cells.Item(132, 10).Value = 2376.0908
$ws.Cells.Item(132, 11).Value = 3135.3489
$ws.Cells.Item(132, 12).Value = 7128.2724
$ws.Cells.Item(132, 13).Value = -605.3489
$ws.Cells.Item(132, 14).Value = -12188.2724

$ws.Cells.Item(136, 8).Value = 1250.579
$ws.Cells.Item(136, 9).Value = 1157.3715
$ws.Cells.Item(136, 10).Value = 2338
$ws.Cells.Item(136, 11).Value = 3472.1145
$ws.Cells.Item(136, 12).Value = 7014
$ws.Cells.Item(136, 13).Value = -922.1144999999997
$ws.Cells.Item(136, 14).Value = -12114

$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()

$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(138, 8).Value = 100000
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 100000
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 100000
$ws.Cells.Item(138, 14).Value = -110280

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 562.5273
$ws.Cells.Item(58, 9).Value = 507.63416
$ws.Cells.Item(58, 10).Value = 723.2857
$ws.Cells.Item(58, 11).Value = 507.63416
$ws.Cells.Item(58, 12).Value = 723.2857
$ws.Cells.Item(58, 13).Value = -304.63416
$ws.Cells.Item(58, 14).Value = -1129.2857

$ws.Cells.Item(132, 8).Value = 2083.0625
$ws.Cells.Item(132, 9).Value = 1637.84
$ws.Cells.Item(132, 10).Value = 3673.1428
$ws.Cells.Item(132, 11).Value = 4913.52
$ws.Cells.Item(132, 12).Value = 11019.4284
$ws.Cells.Item(132, 13).Value = -2383.52
$ws.Cells.Item(132, 14).Value = -16079.4284

$ws.Cells.Item(134, 8).Value = 968.97437
$ws.Cells.Item(134, 9).Value = 971.625
$ws.Cells.Item(134, 10).Value = 956.8570999999999
$ws.Cells.Item(134, 11).Value = 2914.875
$ws.Cells.Item(134, 12).Value = 2870.5713
$ws.Cells.Item(134, 13).Value = -379.875
$ws.Cells.Item(134, 14).Value = -7940.5713

$ws.Cells.Item(136, 8).Value = 562.5273
$ws.Cells.Item(136, 9).Value = 507.63416
$ws.Cells.Item(136, 10).Value = 723.2857
$ws.Cells.Item(136, 11).Value = 1522.90248
$ws.Cells.Item(136, 12).Value = 2169.8571
$ws.Cells.Item(136, 13).Value = 1027.09752
$ws.Cells.Item(136, 14).Value = -7269.8571

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 44.53846
$ws.Cells.Item(38, 9).Value = 45.714287
$ws.Cells.Item(38, 10).Value = 43.166668
$ws.Cells.Item(38, 11).Value = 137.142861
$ws.Cells.Item(38, 12).Value = 129.500004
$ws.Cells.Item(38, 13).Value = 209.857139
$ws.Cells.Item(38, 14).Value = -823.500004

$ws.Cells.Item(86, 8).Value = 561
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 561
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 1683
$ws.Cells.Item(86, 14).Value = -4055

$ws.Cells.Item(89, 8).Value = 561
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 561
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 5049
$ws.Cells.Item(89, 14).Value = -16905

$ws.Cells.Item(121, 8).Value = 832.55554
$ws.Cells.Item(121, 9).Value = 365
$ws.Cells.Item(121, 10).Value = 966.1429000000001
$ws.Cells.Item(121, 11).Value = 1095
$ws.Cells.Item(121, 12).Value = 2898.4287
$ws.Cells.Item(121, 13).Value = 215
$ws.Cells.Item(121, 14).Value = -5518.4287

$ws.Cells.Item(131, 8).Value = 29413072
$ws.Cells.Item(131, 9).Value = 71428780
$ws.Cells.Item(131, 10).Value = 2076.75
$ws.Cells.Item(131, 11).Value = 214286340
$ws.Cells.Item(131, 12).Value = 6230.25
$ws.Cells.Item(131, 13).Value = -214281300
$ws.Cells.Item(131, 14).Value = -16310.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 64289584
$ws.Cells.Item(70, 9).Value = 62504200
$ws.Cells.Item(70, 10).Value = 66670100
$ws.Cells.Item(70, 11).Value = 62504200
$ws.Cells.Item(70, 12).Value = 66670100
$ws.Cells.Item(70, 13).Value = -62503930
$ws.Cells.Item(70, 14).Value = -66670640

$ws.Cells.Item(73, 8).Value = 64289584
$ws.Cells.Item(73, 9).Value = 62504200
$ws.Cells.Item(73, 10).Value = 66670100
$ws.Cells.Item(73, 11).Value = 62504200
$ws.Cells.Item(73, 12).Value = 66670100
$ws.Cells.Item(73, 13).Value = -62503264
$ws.Cells.Item(73, 14).Value = -66671972

$ws.Cells.Item(93, 8).Value = 29999.8
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 29999.8
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 29999.8
$ws.Cells.Item(93, 14).Value = -33743.8

$ws.Cells.Item(132, 8).Value = 1826.4445
$ws.Cells.Item(132, 9).Value = 1365.4615
$ws.Cells.Item(132, 10).Value = 3025
$ws.Cells.Item(132, 11).Value = 4096.3845
$ws.Cells.Item(132, 12).Value = 9075
$ws.Cells.Item(132, 13).Value = -1566.3845
$ws.Cells.Item(132, 14).Value = -14135

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).ClearContents()

$ws.Cells.Item(136, 8).Value = 1496.8096
$ws.Cells.Item(136, 9).Value = 1286.125
$ws.Cells.Item(136, 10).Value = 2171
$ws.Cells.Item(136, 11).Value = 3858.375
$ws.Cells.Item(136, 12).Value = 6513
$ws.Cells.Item(136, 13).Value = -1308.375
$ws.Cells.Item(136, 14).Value = -11613

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 425.25
$ws.Cells.Item(81, 9).Value = 467
$ws.Cells.Item(81, 10).Value = 300
$ws.Cells.Item(81, 11).Value = 934
$ws.Cells.Item(81, 12).Value = 600
$ws.Cells.Item(81, 13).Value = 127
$ws.Cells.Item(81, 14).Value = -2722

$ws.Cells.Item(84, 8).Value = 425.25
$ws.Cells.Item(84, 9).Value = 467
$ws.Cells.Item(84, 10).Value = 300
$ws.Cells.Item(84, 11).Value = 4670
$ws.Cells.Item(84, 12).Value = 3000
$ws.Cells.Item(84, 13).Value = 634
$ws.Cells.Item(84, 14).Value = -13608

$ws.Cells.Item(96, 8).Value = 3522.5715
$ws.Cells.Item(96, 9).Value = 3869.6
$ws.Cells.Item(96, 10).Value = 2655
$ws.Cells.Item(96, 11).Value = 3869.6
$ws.Cells.Item(96, 12).Value = 2655
$ws.Cells.Item(96, 13).Value = -2496.6
$ws.Cells.Item(96, 14).Value = -5401

$ws.Cells.Item(97, 8).Value = 10500
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 10500
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 10500
$ws.Cells.Item(97, 14).Value = -12482

$ws.Cells.Item(132, 8).Value = 1805.1754
$ws.Cells.Item(132, 9).Value = 1798.84
$ws.Cells.Item(132, 10).Value = 1850.4286
$ws.Cells.Item(132, 11).Value = 5396.52
$ws.Cells.Item(132, 12).Value = 5551.2858
$ws.Cells.Item(132, 13).Value = -2866.52
$ws.Cells.Item(132, 14).Value = -10611.2858
